# Refresh of the cryptos price/volume table (GitHub Actions data pull).
# Column D ("Price") holds numeric-looking text (e.g. "116.00", "1.000",
# "26.872.07"); forcing Text format before the write keeps Excel from
# silently coercing it to a Number and dropping significant digits, then
# the style is restored to the sheet default so no formatting leaks in.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.872.07"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.94%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.809.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.95%  "
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.13%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4631"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3702"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  +0.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.71%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.764.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.363"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.02%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.514"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.72"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.07036"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.60%  "
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008687"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.38%  "
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "26.881.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.92%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.314"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.36%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.56"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.40%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.007.93"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.22%  "
$ws.Range("E25").Value = "  -2.78%  "
$ws.Range("E26").Value = "  +0.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.144"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.318"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.46%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "116.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08892"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7510"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.463"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.00%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.920"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.000"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.101"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01965"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.68%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.441"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.69%  "
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05254"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.927"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.65%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5326"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "7.159"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1664"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.484"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.23%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4964"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.38%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.36%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.673"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.71%  "
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.000"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.16%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "103.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06288"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.43%  "
